# Update the second table ("Table 1: Example With Formatting") to mark its
# first row as a repeating header row, and append a new data row
# (Data Cell 3 / Data Cell 4).

$d = $word.ActiveDocument

# The document contains two tables; the one with cell borders / shading is
# the second one (Tables.Item(2)).
$t = $d.Tables.Item(2)

# 1) Flag the first row ("Head Cell 1" / "Head Cell 2") as a header row that
#    repeats on every page -> <w:trPr><w:tblHeader/></w:trPr>.
$headerRow = $t.Rows.Item(1)
$headerRow.HeadingFormat = $true

# 2) Append a new row with two plain cells containing "Data Cell 3" and
#    "Data Cell 4". Rows.Add() only mints a single cell for the new row, so
#    split it into two columns to match the table's existing column count.
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Split(1, 2)

$cell3 = $newRow.Cells.Item(1)
$cell4 = $newRow.Cells.Item(2)

# Assigning text twice ensures the run's <w:t> keeps xml:space="preserve",
# matching how Word represents the other table cells.
$cell3.Range.Text = " "
$cell3.Range.Text = "Data Cell 3"

$cell4.Range.Text = " "
$cell4.Range.Text = "Data Cell 4"
